$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL: http://fhir.meuhedet.co.il/code/lab-result -> new IG url
$ws.Range("B2").Value = "http://ig.fhir-il-community.org/LRC/CodeSystem/meuhedet-lab-result-code"

# Date: 2024-03-03T10:40:33+02:00 -> 2024-03-04T15:01:53+02:00
$ws.Range("B8").Value = "2024-03-04T15:01:53+02:00"

# Publisher: Outburn LTD. -> FHIR-il community
$ws.Range("B9").Value = "FHIR-il community"

# Content: fragment -> complete
$ws.Range("B19").Value = "complete"
